$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 9, 10 and 12 (row 11 stays put):
#   new row 9  <- old row 12
#   new row 10 <- old row 9
#   new row 12 <- old row 10
# Capture the "before" values for the columns that actually change,
# then write them back in the rotated order.

$numericCols = @("A","B","E","Q","R")
$textCols = @("F","G","H","I","Y","AA","AC","AX")
$cols = $numericCols + $textCols

$row9 = @{}
$row10 = @{}
$row12 = @{}

foreach ($col in $cols) {
    $row9[$col] = $ws.Range("$col`9").Value()
    $row10[$col] = $ws.Range("$col`10").Value()
    $row12[$col] = $ws.Range("$col`12").Value()
}

function Set-RotatedValue($col, $rowNum, $value) {
    $cell = $ws.Range("$col$rowNum")
    if ($textCols -contains $col) {
        # Force literal text so Excel doesn't auto-convert look-alike
        # numbers/dates (e.g. "3" or "2023-08-12") to numeric/date types.
        $cell.NumberFormat = "@"
        $cell.Value = [string]$value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

foreach ($col in $cols) {
    Set-RotatedValue $col 9 $row12[$col]
    Set-RotatedValue $col 10 $row9[$col]
    Set-RotatedValue $col 12 $row10[$col]
}
